$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top; this shifts all existing rows down by one.
$ws.Rows.Item(1).Insert()

# New header row
$ws.Range("A1").Value = "WebElementName"
$ws.Range("B1").Value = "Name"

# New values inserted into existing (now shifted) rows
$ws.Range("B2").Value = "abc"
$ws.Range("B3").Value = "efd"
